$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns stay plain text so numeric-looking strings
# (e.g. "1.010", "30.541.92") are not reinterpreted as numbers/dates.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.547.32"
$ws.Range("E2").Value = "  +2.74%  "
$ws.Range("D3").Value = "2.128.50"
$ws.Range("E3").Value = "  +1.64%  "
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "347.74"
$ws.Range("E5").Value = "  +0.72%  "
$ws.Range("D6").Value = "1.008"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "0.5265"
$ws.Range("E7").Value = "  +1.82%  "
$ws.Range("D8").Value = "0.4495"
$ws.Range("E8").Value = "  +1.31%  "
$ws.Range("D9").Value = "54.19"
$ws.Range("E9").Value = "  +5.08%  "
$ws.Range("D10").Value = "0.09414"
$ws.Range("E10").Value = "  -0.71%  "
$ws.Range("D11").Value = "1.188"
$ws.Range("E11").Value = "  +0.82%  "
$ws.Range("D12").Value = "25.36"
$ws.Range("E12").Value = "  -0.72%  "
$ws.Range("D13").Value = "8.793"
$ws.Range("E13").Value = "  +8.54%  "
$ws.Range("D14").Value = "7.003"
$ws.Range("E14").Value = "  +3.39%  "
$ws.Range("D15").Value = "2.098.26"
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("D16").Value = "102.59"
$ws.Range("E16").Value = "  +3.10%  "
$ws.Range("D17").Value = "0.00001174"
$ws.Range("E17").Value = "  +0.47%  "
$ws.Range("D18").Value = "1.009"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("E19").Value = "  +2.64%  "
$ws.Range("D20").Value = "0.06738"
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("D21").Value = "6.363"
$ws.Range("E21").Value = "  +2.54%  "
$ws.Range("D22").Value = "1.008"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").Value = "30.573.84"
$ws.Range("E23").Value = "  +2.59%  "
$ws.Range("D24").Value = "12.79"
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("D25").Value = "2.338"
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("D26").Value = "2.374.96"
$ws.Range("E26").Value = "  +1.34%  "
$ws.Range("D27").Value = "22.29"
$ws.Range("E27").Value = "  +1.54%  "
$ws.Range("D29").Value = "163.36"
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").Value = "135.03"
$ws.Range("E30").Value = "  +1.50%  "
$ws.Range("D31").Value = "1.166"
$ws.Range("E31").Value = "  +1.37%  "
$ws.Range("D32").Value = "1.783"
$ws.Range("E32").Value = "  +9.98%  "
$ws.Range("D33").Value = "0.1065"
$ws.Range("E33").Value = "  +0.95%  "
$ws.Range("D34").Value = "6.970"
$ws.Range("E34").Value = "  +12.86%  "
$ws.Range("D35").Value = "6.322"
$ws.Range("E35").Value = "  +1.74%  "
$ws.Range("D37").Value = "10.72"
$ws.Range("E37").Value = "  +5.76%  "
$ws.Range("D38").Value = "0.02665"
$ws.Range("E38").Value = "  +3.39%  "
$ws.Range("D39").Value = "0.06888"
$ws.Range("E39").Value = "  +2.21%  "
$ws.Range("D40").Value = "0.7157"
$ws.Range("E40").Value = "  +3.90%  "
$ws.Range("D41").Value = "12.74"
$ws.Range("E41").Value = "  +2.65%  "
$ws.Range("D42").Value = "0.2260"
$ws.Range("E42").Value = "  -0.84%  "
$ws.Range("D43").Value = "1.338"
$ws.Range("E43").Value = "  +4.29%  "
$ws.Range("D44").Value = "0.6951"
$ws.Range("E44").Value = "  +4.20%  "
$ws.Range("D45").Value = "14.72"
$ws.Range("E45").Value = "  +3.92%  "
$ws.Range("D46").Value = "2.410"
$ws.Range("E46").Value = "  +4.61%  "
$ws.Range("D47").Value = "1.008"
$ws.Range("E47").Value = "  +0.15%  "
$ws.Range("E48").Value = "  +11.76%  "
$ws.Range("D49").Value = "3.652"
$ws.Range("E49").Value = "  +0.66%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Value = "1.233"
$ws.Range("E50").Value = "  +1.19%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.00000000345"
$ws.Range("E51").Value = "  +0.15%  "
